$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (shifts existing A:I to C:K)
$ws.Range("A:B").Insert()

# New header cells
$ws.Range("A1").Value = "Emp_Id"
$ws.Range("B1").Value = "FULL_NAME"

# Fill Emp_Id (1..12) and FULL_NAME (=FIRST_NAME & " " & LAST_NAME) for rows 2..13
for ($i = 2; $i -le 13; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $firstName = $ws.Cells.Item($i, 3).Value()
    $lastName = $ws.Cells.Item($i, 4).Value()
    $ws.Cells.Item($i, 2).Value = "$firstName $lastName"
}

# AutoFit the new columns to match bestFit widths
$ws.Range("A:B").Columns.AutoFit()

# Update selection to match target state (whole column F selected)
$ws.Range("F1:F1048576").Select()
